# Updated cryptos list on Wed Nov 29 19:40:11 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for the crypto table on
# the active sheet, and swaps the Monero/Cosmos rows (26/27) to reflect
# their new ranking order. Cells whose new text looks like a plain number
# ("227.30", "0.613", ...) are written with a leading apostrophe so Excel
# keeps storing them as text (matching the original inlineStr cells)
# instead of silently converting them to numeric values; the apostrophe
# marker is then cleared by resetting the cell style back to "Normal" so
# no stray number-format/quote-prefix style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.861.01'
$ws.Range('D3').Value = '2.035.31'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'227.30"
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').Value = "'0.613"
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').Value = "'60.08"
$ws.Range('E7').Value = '  +3.32%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').Value = "'0.0818"
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = "'14.62"
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '2.337.00'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').Value = "'21.01"
$ws.Range('E14').Value = '  +1.42%  '
$ws.Range('D15').Value = "'0.759"
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '2.037.92'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = '37.832.96'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = "'6.06"
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('D20').Value = "'69.77"
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '0.0₃0822'
$ws.Range('E21').Value = '  -1.18%  '
$ws.Range('D22').Value = "'224.47"
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = "'2.37"
$ws.Range('E24').Value = '  -3.34%  '
$ws.Range('D25').Value = "'2.20"
$ws.Range('E25').Value = '  -2.01%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = "'9.26"
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'165.10"
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = "'0.129"
$ws.Range('E28').Value = '  -3.96%  '
$ws.Range('D29').Value = "'18.87"
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D31').Value = "'0.119"
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('E32').Value = '  -2.70%  '
$ws.Range('E33').Value = '  +4.41%  '
$ws.Range('E34').Value = '  -2.63%  '
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').Value = "'6.41"
$ws.Range('E36').Value = '  +5.07%  '
$ws.Range('D37').Value = "'2.25"
$ws.Range('E37').Value = '  -5.66%  '
$ws.Range('E38').Value = '  -2.76%  '
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').Value = '1.540.69'
$ws.Range('E40').Value = '  +3.74%  '
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('D42').Value = "'96.94"
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('D43').Value = "'16.86"
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').Value = "'0.0922"
$ws.Range('E45').Value = '  -2.12%  '
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').Value = "'3.92"
$ws.Range('E47').Value = '  -4.38%  '
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').Value = '2.227.56'
$ws.Range('E51').Value = '  -1.05%  '

# Reset style on cells that required a quote-prefix to stay text,
# so no stray style index is left behind (matches the source diff,
# which only touches cell text, not styles).
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
